$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Testi 7.4.2017 @ 21:47"
